$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates derived from the latest cryptos-list refresh.
# Only cells that actually changed are listed per row; D (Price) keeps its
# original plain-text formatting (so "1.0000" etc. do not get coerced to numbers).
$updates = @(
    @{ Row=2; D="31.247.22"; E="  +1.51%  " }
    @{ Row=3; D="1.965.22"; E="  +0.80%  " }
    @{ Row=4; D="1.001"; E="  +0.12%  " }
    @{ Row=5; D="246.12"; E="  -0.80%  " }
    @{ Row=6; D="1.0000"; E="  -0.01%  " }
    @{ Row=7; D="0.4923"; E="  +2.36%  " }
    @{ Row=8; D="44.86"; E="  -1.57%  " }
    @{ Row=9; D="0.2992"; E="  +1.27%  " }
    @{ Row=10; D="0.06890"; E="  +0.92%  " }
    @{ Row=11; E="  -0.91%  " }
    @{ Row=12; D="109.21"; E="  -3.33%  " }
    @{ Row=13; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.07763"; E="  +1.38%  " }
    @{ Row=14; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.929.77"; E="  -1.06%  " }
    @{ Row=15; D="5.486"; E="  -1.51%  " }
    @{ Row=16; D="0.7184"; E="  +3.56%  " }
    @{ Row=17; D="288.07"; E="  -3.61%  " }
    @{ Row=18; D="31.152.28"; E="  +1.33%  " }
    @{ Row=19; D="0.000007816"; E="  +1.45%  " }
    @{ Row=20; D="13.32"; E="  +0.08%  " }
    @{ Row=21; D="2.201.35"; E="  -0.09%  " }
    @{ Row=22; D="1.000"; E="  +0.01%  " }
    @{ Row=23; D="5.548"; E="  -2.67%  " }
    @{ Row=24; D="1.001"; E="  +0.17%  " }
    @{ Row=25; D="6.603"; E="  +0.16%  " }
    @{ Row=26; D="9.912"; E="  +1.65%  " }
    @{ Row=27; D="169.73"; E="  +1.02%  " }
    @{ Row=28; D="20.39"; E="  -0.46%  " }
    @{ Row=29; D="2.226"; E="  +1.95%  " }
    @{ Row=30; D="0.1059"; E="  -3.02%  " }
    @{ Row=31; D="1.433"; E="  +0.11%  " }
    @{ Row=32; D="1.585"; E="  -0.73%  " }
    @{ Row=33; D="4.647"; E="  +0.35%  " }
    @{ Row=34; D="4.481"; E="  +0.78%  " }
    @{ Row=35; D="0.04994"; E="  -1.65%  " }
    @{ Row=36; D="0.7664"; E="  -1.75%  " }
    @{ Row=37; D="1.187"; E="  +1.73%  " }
    @{ Row=38; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02053"; E="  -1.12%  " }
    @{ Row=39; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.731"; E="  -0.03%  " }
    @{ Row=40; D="2.709"; E="  +0.08%  " }
    @{ Row=41; D="2.201"; E="  +7.59%  " }
    @{ Row=42; D="6.446"; E="  +7.93%  " }
    @{ Row=43; D="0.4564"; E="  +1.80%  " }
    @{ Row=44; D="109.86"; E="  -1.41%  " }
    @{ Row=45; D="0.8855"; E="  +0.91%  " }
    @{ Row=46; D="72.49"; E="  +1.04%  " }
    @{ Row=47; D="8.157"; E="  +9.81%  " }
    @{ Row=48; D="1.001"; E="  -0.22%  " }
    @{ Row=49; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="969.67"; E="  +6.37%  " }
    @{ Row=50; E="  +0.90%  " }
    @{ Row=51; D="0.2646"; E="  +3.21%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Force text format so numeric-looking strings are not auto-converted
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
